$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.165.82'
$ws.Range("E2").Value = '  +3.02%  '

# Row 3
$ws.Range("D3").Value = '3.625.69'
$ws.Range("E3").Value = '  +2.43%  '

# Row 4
$ws.Range("E4").Value = '  +0.22%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '624.17'
$ws.Range("E5").Value = '  +2.53%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.10'
$ws.Range("E6").Value = '  +2.46%  '

# Row 7
$ws.Range("D7").Value = '3.627.20'
$ws.Range("E7").Value = '  +2.55%  '

# Row 8
$ws.Range("E8").Value = '  -0.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.491'
$ws.Range("E9").Value = '  +1.44%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").Value = '  +1.97%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.16'
$ws.Range("E11").Value = '  +4.08%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.437'
$ws.Range("E12").Value = '  +1.87%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000223'
$ws.Range("E13").Value = '  +0.97%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.12'
$ws.Range("E14").Value = '  +3.79%  '

# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '4.240.52'
$ws.Range("E15").Value = '  +2.51%  '

# Row 16
$ws.Range("D16").Value = '3.615.77'
$ws.Range("E16").Value = '  +2.07%  '

# Row 17
$ws.Range("D17").Value = '69.413.53'
$ws.Range("E17").Value = '  +3.55%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.117'
$ws.Range("E18").Value = '  -0.37%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.57'
$ws.Range("E19").Value = '  +3.89%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.74'
$ws.Range("E20").Value = '  +2.02%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.18'
$ws.Range("E21").Value = '  +9.66%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '464.27'
$ws.Range("E22").Value = '  +3.66%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.639'
$ws.Range("E23").Value = '  +0.91%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.14'
$ws.Range("E24").Value = '  -0.20%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000135'
$ws.Range("E25").Value = '  +10.01%  '

# Row 26
$ws.Range("D26").Value = '3.766.82'
$ws.Range("E26").Value = '  +2.39%  '

# Row 27
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.45'
$ws.Range("E27").Value = '  +1.68%  '

# Row 28
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.11'
$ws.Range("E29").Value = '  +10.59%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.60'
$ws.Range("E30").Value = '  +2.44%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.72'
$ws.Range("E31").Value = '  +3.09%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.173'
$ws.Range("E32").Value = '  +9.55%  '

# Row 33
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.53'
$ws.Range("E33").Value = '  +5.79%  '

# Row 34
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  -0.35%  '

# Row 35
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.45'
$ws.Range("E35").Value = '  +2.67%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.94'
$ws.Range("E36").Value = '  +2.35%  '

# Row 37
$ws.Range("B37").Value = 'RenzoRestakedETH'
$ws.Range("C37").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D37").Value = '3.609.99'
$ws.Range("E37").Value = '  +2.18%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.30'
$ws.Range("E38").Value = '  +3.31%  '

# Row 39
$ws.Range("E39").Value = '  +0.03%  '

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.33'
$ws.Range("E40").Value = '  +8.49%  '

# Row 41
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0926'
$ws.Range("E41").Value = '  +6.70%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.51%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '174.94'
$ws.Range("E43").Value = '  -0.62%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.61'
$ws.Range("E44").Value = '  +0.62%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.914'
$ws.Range("E45").Value = '  +2.43%  '

# Row 46
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '30.49'
$ws.Range("E46").Value = '  +9.07%  '

# Row 47
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.36'
$ws.Range("E47").Value = '  +10.57%  '

# Row 48
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.10'
$ws.Range("E48").Value = '  +0.71%  '

# Row 49
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.76'
$ws.Range("E49").Value = '  +4.84%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.75'
$ws.Range("E50").Value = '  +1.94%  '

# Row 51
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.266'
$ws.Range("E51").Value = '  +6.38%  '

